$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction on SMOTE_Aug_cc32 (IDs 10 and 11 swapped):
# the "Meaning" text for AE_N / AI_N corner-case rows referenced the wrong
# cc_xx naming (cc_N0 / cc_X0) -- should be cc_0N / cc_0X.
$b36 = $ws.Range("B36").Value2
$ws.Range("B36").Value = $b36.Replace("cc_N0", "cc_0N")

$b37 = $ws.Range("B37").Value2
$ws.Range("B37").Value = $b37.Replace("cc_X0", "cc_0X")

# Reflect the author's final on-screen selection (the merged title cell).
$ws.Range("A1:C1").Select()
